$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estoque")

# Update quantities and alteration dates for a few existing rows
$ws.Range("C5").Value = 32
$ws.Range("E5").Value = "21/10/2025 17:57"

$ws.Range("C7").Value = 35
$ws.Range("E7").Value = "21/10/2025 19:51"

$ws.Range("C15").Value = 16
$ws.Range("E15").Value = "21/10/2025 19:58"

# Remove the now-obsolete duplicate row (previously row 16)
$ws.Rows.Item(16).Delete()
